{"js": "// FIX: \"at age\" and \"at length\" were the wrong way around.\n//\n// The document is a template describing the \"proportions_tag_at_age\"\n// observation/process. Every occurrence of the word \"length\" in the body\n// text is really meant to say \"age\" (title, the \"proportions_at_length\"\n// identifier, the \"tag_by_length\" identifier, etc.) \u2014 the commit simply\n// swaps those occurrences of \"length\" for \"age\" while leaving every other\n// character (capitalisation, punctuation, surrounding text, formatting)\n// untouched.\n\nconst body = context.document.body;\n\n// Find every run of the literal text \"length\" anywhere in the document\n// body (titles, paragraphs, etc.) and swap it for \"age\" in place so the\n// surrounding formatting (font size, bold, colour, proofing tags, ...)\n// is preserved.\nconst results = body.search(\"length\", { matchCase: false, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"age\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# FIX: \"at age\" and \"at length\" were the wrong way around.\n#\n# The document is a template describing the \"proportions_tag_at_age\"\n# observation/process. Every occurrence of the word \"length\" in the body\n# text is really meant to say \"age\" (the big title, the\n# \"proportions_at_length\" identifier, the \"tag_by_length\" identifier,\n# etc.) -- the commit simply swaps those occurrences of \"length\" for\n# \"age\" while leaving every other character (capitalisation, punctuation,\n# surrounding text, formatting) untouched.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"length\"\n$find.Replacement.Text = \"age\"\n$find.MatchCase = $false\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
